# Update the 2025Q3 row (row 29) metrics in the quarterly recurrence sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C29").Value = 195
$ws.Range("D29").Value = 29
$ws.Range("E29").Value = 166
$ws.Range("F29").Value = 4.991394148020654
